# Regenerate merged AHB files
# 1) Rename the FV-specific header columns (old -> FV2310, new -> FV2404)
# 2) Freeze the header row (pane split under row 1)
# 3) Turn the A1:U94 range into an Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header renames -----------------------------------------------------
$baseCols = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

$headers = New-Object System.Collections.ArrayList
foreach ($col in $baseCols) { [void]$headers.Add($col + "_FV2310") }
[void]$headers.Add("diff")
foreach ($col in $baseCols) { [void]$headers.Add($col + "_FV2404") }

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2) Create the table BEFORE re-bolding the header row ------------------
# (the header row already carries a bold style from the template; creating
# the table while it is still bold makes Excel capture that as an explicit
# header dxf, which the source workbook does not have, so we temporarily
# clear it, build the table, then restore the bold header formatting)
$headerRange = $ws.Range("A1:U1")
$wasBold = $headerRange.Font.Bold
$headerRange.Font.Bold = $false

$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U94"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

$headerRange.Font.Bold = $wasBold

# --- 3) Freeze the header row ----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
